$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new feed rows were picked up by the filtering workflow: the same
# Geneseeq/Roche Rozlytrek story as reported by GenomeWeb and 360Dx.
$title = "Geneseeq Nabs China NMPA Approval for Pan-Solid Tumor Test as CDx for Roche's Rozlytrek"

$urlGenomeweb = "https://www.genomeweb.com/cancer/geneseeq-nabs-china-nmpa-approval-pan-solid-tumor-test-cdx-roches-rozlytrek"
$url360dx = "https://www.360dx.com/cancer/geneseeq-nabs-china-nmpa-approval-pan-solid-tumor-test-cdx-roches-rozlytrek"

# Row 59: GenomeWeb link
$ws.Cells.Item(59, 2).Value = "CDx"
$ws.Cells.Item(59, 3).Value = $title
$ws.Hyperlinks.Add($ws.Cells.Item(59, 1), $urlGenomeweb, "", "", $urlGenomeweb)
$ws.Cells.Item(59, 1).Style = "Hyperlink"

# Row 60: 360Dx link (same keywords/title)
$ws.Cells.Item(60, 2).Value = "CDx"
$ws.Cells.Item(60, 3).Value = $title
$ws.Hyperlinks.Add($ws.Cells.Item(60, 1), $url360dx, "", "", $url360dx)
$ws.Cells.Item(60, 1).Style = "Hyperlink"

Write-Host "Added 2 filtered-feed rows"
